$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.260.14'
$ws.Range("E2").Value = '  -1.41%  '
$ws.Range("D3").Value = '2.426.18'
$ws.Range("E3").Value = '  -0.87%  '
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").Value = '''570.46'
$ws.Range("E5").Value = '  -2.06%  '
$ws.Range("D6").Value = '''140.12'
$ws.Range("E6").Value = '  -1.94%  '
$ws.Range("E7").Value = '  +0.17%  '
$ws.Range("E8").Value = '  -0.90%  '
$ws.Range("D9").Value = '2.412.07'
$ws.Range("E9").Value = '  -1.25%  '
$ws.Range("E11").Value = '  -0.24%  '
$ws.Range("E12").Value = '  -2.69%  '
$ws.Range("E13").Value = '  -0.88%  '
$ws.Range("D14").Value = '''26.19'
$ws.Range("E14").Value = '  -0.90%  '
$ws.Range("D15").Value = '''0.0000172'
$ws.Range("E15").Value = '  -2.20%  '
$ws.Range("D16").Value = '2.841.40'
$ws.Range("D17").Value = '61.142.53'
$ws.Range("E17").Value = '  -1.45%  '
$ws.Range("D18").Value = '2.411.36'
$ws.Range("E18").Value = '  -0.91%  '
$ws.Range("D19").Value = '''7.78'
$ws.Range("E19").Value = '  +7.65%  '
$ws.Range("D20").Value = '''10.66'
$ws.Range("E20").Value = '  -0.36%  '
$ws.Range("D21").Value = '''323.83'
$ws.Range("E21").Value = '  -0.73%  '
$ws.Range("D22").Value = '''4.07'
$ws.Range("E22").Value = '  -0.82%  '
$ws.Range("D23").Value = '''6.10'
$ws.Range("E23").Value = '  +2.13%  '
$ws.Range("E24").Value = '  +0.07%  '
$ws.Range("E25").Value = '  -2.99%  '
$ws.Range("D26").Value = '''64.69'
$ws.Range("E26").Value = '  -1.39%  '
$ws.Range("D27").Value = '''592.76'
$ws.Range("E27").Value = '  -0.99%  '
$ws.Range("D28").Value = '''8.28'
$ws.Range("E28").Value = '  -9.18%  '
$ws.Range("E29").Value = '  -0.94%  '
$ws.Range("D30").Value = '0.0₃0939'
$ws.Range("E30").Value = '  -2.92%  '
$ws.Range("D31").Value = '''7.94'
$ws.Range("E31").Value = '  -0.65%  '
$ws.Range("E32").Value = '  -4.34%  '
$ws.Range("E33").Value = '  -4.02%  '
$ws.Range("E34").Value = '  -0.91%  '
$ws.Range("E35").Value = '  -0.07%  '
$ws.Range("D36").Value = '''1.43'
$ws.Range("D37").Value = '''4.63'
$ws.Range("E37").Value = '  -5.04%  '
$ws.Range("D38").Value = '''151.91'
$ws.Range("E38").Value = '  -0.59%  '
$ws.Range("D39").Value = '''0.369'
$ws.Range("D40").Value = '''18.25'
$ws.Range("E40").Value = '  -0.89%  '
$ws.Range("D41").Value = '''5.19'
$ws.Range("E41").Value = '  -1.99%  '
$ws.Range("E42").Value = '  +0.02%  '
$ws.Range("E43").Value = '  -1.76%  '
$ws.Range("D44").Value = '''41.24'
$ws.Range("E44").Value = '  -4.50%  '
$ws.Range("D45").Value = '''2.38'
$ws.Range("E45").Value = '  -5.23%  '
$ws.Range("D46").Value = '0.0₆0298'
$ws.Range("E46").Value = '  +7.77%  '
$ws.Range("D47").Value = '''143.37'
$ws.Range("E47").Value = '  +1.17%  '
$ws.Range("E48").Value = '  -2.46%  '
$ws.Range("D49").Value = '''0.589'
$ws.Range("E49").Value = '  -2.14%  '
$ws.Range("D50").Value = '''19.60'
$ws.Range("E50").Value = '  -1.28%  '
$ws.Range("E51").Value = '  -3.00%  '
